$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('I1').Value = 'Other found locations'

$ws.Range('E2').Value = '[Alfredo%Guiroy%NULL%1,     Martín%Gagliardi%NULL%2,     Martín%Gagliardi%NULL%0,     Nicolas%Coombes%NULL%1,     Federico%Landriel%NULL%2,     Federico%Landriel%NULL%0,     Carlos%Zanardi%NULL%1,     Gastón Camino%Willhuber%NULL%1,     Juan Pablo%Guyot%NULL%1,     Marcelo%Valacco%NULL%1]'
$ws.Range('I2').Value = '_PMC'

$ws.Range('E3').Value = '[Yufei%Li%NULL%10,     Nathaniel%Scherer%NULL%10,     Lambert%Felix%NULL%20,     Lambert%Felix%NULL%0,     Hannah%Kuper%NULL%20,     Hannah%Kuper%NULL%0,     Jakob%Pietschnig%NULL%20,     Jakob%Pietschnig%NULL%0]'
$ws.Range('I3').Value = '_PMC'

$ws.Range('I4').Value = ''

$ws.Range('E5').Value = '[Yufei%Li%NULL%0,     Nathaniel%Scherer%NULL%0,     Lambert%Felix%NULL%0,     Lambert%Felix%NULL%0,     Hannah%Kuper%NULL%0,     Hannah%Kuper%NULL%0,     Jakob%Pietschnig%NULL%0,     Jakob%Pietschnig%NULL%0]'
$ws.Range('I5').Value = '_PMC'

$ws.Range('E6').Value = '[Kossi Blewussi%Kounou%benkounou@hotmail.fr%1,     Koffi Mawuse%Guédénon%NULL%2,     Koffi Mawuse%Guédénon%NULL%0,     Ayoko Akouavi%Dogbe Foli%NULL%1,     Eric%Gnassounou‐Akpa%NULL%1]'
$ws.Range('I6').Value = '_PMC'

$ws.Range('E7').Value = '[Jianbo%Lai%NULL%0,     Simeng%Ma%NULL%1,     Ying%Wang%NULL%3,     Zhongxiang%Cai%NULL%1,     Jianbo%Hu%NULL%1,     Ning%Wei%NULL%1,     Jiang%Wu%NULL%1,     Hui%Du%NULL%1,     Tingting%Chen%NULL%1,     Ruiting%Li%NULL%1,     Huawei%Tan%NULL%1,     Lijun%Kang%NULL%1,     Lihua%Yao%NULL%1,     Manli%Huang%NULL%1,     Huafen%Wang%NULL%1,     Gaohua%Wang%NULL%1,     Zhongchun%Liu%NULL%1,     Shaohua%Hu%NULL%1]'
$ws.Range('I7').Value = '_PMC'

$ws.Range('E8').Value = '[Simon Ching%Lam%NULL%1,     Teresa%Arora%NULL%1,     Ian%Grey%NULL%1,     Lorna Kwai Ping%Suen%NULL%1,     Emma Yun-zhi%Huang%NULL%1,     Daofan%Li%NULL%1,     Kin Bong Hubert%Lam%NULL%1]'
$ws.Range('I8').Value = '_PMC'

$ws.Range('E9').Value = '[Yufei%Li%NULL%0,     Nathaniel%Scherer%NULL%0,     Lambert%Felix%NULL%0,     Lambert%Felix%NULL%0,     Hannah%Kuper%NULL%0,     Hannah%Kuper%NULL%0,     Jakob%Pietschnig%NULL%0,     Jakob%Pietschnig%NULL%0]'
$ws.Range('I9').Value = '_PMC'

$ws.Range('E10').Value = '[Ruilin%Li%NULL%1,     Youlin%Chen%NULL%1,     Jianlin%Lv%NULL%1,     Linlin%Liu%NULL%1,     Shiqin%Zong%NULL%1,     Hanxia%Li%NULL%1,     Hong%Li%NULL%1,     Massimo%Tusconi.%NULL%2,     Massimo%Tusconi.%NULL%0]'
$ws.Range('I10').Value = '_PMC'

$ws.Range('E11').Value = '[Yufei%Li%NULL%0,     Nathaniel%Scherer%NULL%0,     Lambert%Felix%NULL%0,     Lambert%Felix%NULL%0,     Hannah%Kuper%NULL%0,     Hannah%Kuper%NULL%0,     Jakob%Pietschnig%NULL%0,     Jakob%Pietschnig%NULL%0]'
$ws.Range('I11').Value = '_PMC'

$ws.Range('E12').Value = '[Wen%Lu%luwen67@sina.com%0,     Hang%Wang%NULL%1,     Yuxing%Lin%NULL%1,     Li%Li%lilifuzhou@126.com%0]'
$ws.Range('I12').Value = '_PMC_elsevier'

$ws.Range('E13').Value = '[Abdallah Y.%Naser%abdallah.naser@iu.edu.jo%0,     Eman Zmaily%Dahmash%NULL%2,     Eman Zmaily%Dahmash%NULL%0,     Rabaa%Al‐Rousan%NULL%2,     Rabaa%Al‐Rousan%NULL%0,     Hassan%Alwafi%NULL%1,     Hamzeh Mohammad%Alrawashdeh%NULL%1,     Imene%Ghoul%NULL%1,     Anwer%Abidine%NULL%1,     Mohammed A.%Bokhary%NULL%1,     Hadeel T.%AL‐Hadithi%NULL%1,     Dalia%Ali%NULL%1,     Rasha%Abuthawabeh%NULL%1,     Ghada Mohammad%Abdelwahab%NULL%1,     Yosra J.%Alhartani%NULL%1,     Haneen%Al Muhaisen%NULL%1,     Ayah%Dagash%NULL%1,     Hamad S.%Alyami%NULL%1]'
$ws.Range('I13').Value = '_PMC'

$ws.Range('E14').Value = '[Moluk%Pouralizadeh%NULL%1,     Zahra%Bostani%NULL%1,     Saman%Maroufizadeh%NULL%1,     Atefeh%Ghanbari%NULL%1,     Maryam%Khoshbakht%NULL%1,     Seyed Amirhossein%Alavi%NULL%1,     Sadra%Ashrafi%NULL%1]'
$ws.Range('I14').Value = '_PMC_elsevier'

$ws.Range('I15').Value = ''

$ws.Range('E16').Value = '[Jianyu%Que%NULL%1,     Le%Shi%NULL%2,     Le%Shi%NULL%0,     Jiahui%Deng%NULL%1,     Jiajia%Liu%NULL%1,     Li%Zhang%NULL%0,     Suying%Wu%NULL%1,     Yimiao%Gong%NULL%1,     Weizhen%Huang%NULL%1,     Kai%Yuan%NULL%1,     Wei%Yan%NULL%1,     Yankun%Sun%NULL%1,     Maosheng%Ran%NULL%1,     Yanping%Bao%NULL%1,     Lin%Lu%NULL%1]'
$ws.Range('I16').Value = '_PMC'

$ws.Range('E17').Value = '[Rodolfo%Rossi%NULL%0,     Valentina%Socci%NULL%1,     Francesca%Pacitti%NULL%1,     Giorgio%Di Lorenzo%NULL%1,     Antinisca%Di Marco%NULL%1,     Alberto%Siracusano%NULL%1,     Alessandro%Rossi%NULL%1]'
$ws.Range('I17').Value = '_PMC'

$ws.Range('E18').Value = '[ M.%Salman%null%1,      M. H.% Raza%null%1,      Z. U.% Mustafa%null%1,      T. M.% Khan%null%1,      N.% Asif%null%1,      H.% Tahir%null%1,      N.% Shehzadi%null%1,      K. % Hussain%null%1]'
$ws.Range('F18').Value = 'not found'
$ws.Range('G18').Value = 'N/A'
$ws.Range('I18').Value = '_MedBiorxiv'

$ws.Range('E19').Value = '[Ari%Shechter%NULL%1,     Franchesca%Diaz%NULL%1,     Nathalie%Moise%NULL%1,     D. Edmund%Anstey%NULL%1,     Siqin%Ye%NULL%1,     Sachin%Agarwal%NULL%1,     Jeffrey L.%Birk%NULL%1,     Daniel%Brodie%NULL%3,     Diane E.%Cannone%NULL%1,     Bernard%Chang%NULL%1,     Jan%Claassen%NULL%2,     Talea%Cornelius%NULL%1,     Lilly%Derby%NULL%1,     Melissa%Dong%NULL%1,     Raymond C.%Givens%NULL%1,     Beth%Hochman%NULL%1,     Shunichi%Homma%NULL%1,     Ian M.%Kronish%NULL%1,     Sung A.J.%Lee%NULL%1,     Wilhelmina%Manzano%NULL%1,     Laurel E.S.%Mayer%NULL%1,     Cara L.%McMurry%NULL%1,     Vivek%Moitra%NULL%1,     Patrick%Pham%NULL%1,     LeRoy%Rabbani%NULL%1,     Reynaldo R.%Rivera%NULL%1,     Allan%Schwartz%NULL%2,     Joseph E.%Schwartz%NULL%1,     Peter A.%Shapiro%NULL%1,     Kaitlin%Shaw%NULL%1,     Alexandra M.%Sullivan%NULL%1,     Courtney%Vose%NULL%1,     Lauren%Wasson%NULL%1,     Donald%Edmondson%NULL%1,     Marwah%Abdalla%NULL%1]'
$ws.Range('I19').Value = '_PMC_elsevier'

$ws.Range('E20').Value = '[Xingyue%Song%NULL%0,     Wenning%Fu%NULL%1,     Xiaoran%Liu%NULL%1,     Zhiqian%Luo%NULL%1,     Rixing%Wang%NULL%1,     Ning%Zhou%NULL%1,     Shijiao%Yan%NULL%2,     Chuanzhu%Lv%NULL%2]'
$ws.Range('I20').Value = '_PMC_elsevier'

$ws.Range('E21').Value = '[Yufei%Li%NULL%0,     Nathaniel%Scherer%NULL%0,     Lambert%Felix%NULL%0,     Lambert%Felix%NULL%0,     Hannah%Kuper%NULL%0,     Hannah%Kuper%NULL%0,     Jakob%Pietschnig%NULL%0,     Jakob%Pietschnig%NULL%0]'
$ws.Range('I21').Value = '_PMC'

$ws.Range('E22').Value = '[Yufei%Li%NULL%0,     Nathaniel%Scherer%NULL%0,     Lambert%Felix%NULL%0,     Lambert%Felix%NULL%0,     Hannah%Kuper%NULL%0,     Hannah%Kuper%NULL%0,     Jakob%Pietschnig%NULL%0,     Jakob%Pietschnig%NULL%0]'
$ws.Range('I22').Value = '_PMC'

$ws.Range('E23').Value = '[Yufei%Li%NULL%0,     Nathaniel%Scherer%NULL%0,     Lambert%Felix%NULL%0,     Lambert%Felix%NULL%0,     Hannah%Kuper%NULL%0,     Hannah%Kuper%NULL%0,     Jakob%Pietschnig%NULL%0,     Jakob%Pietschnig%NULL%0]'
$ws.Range('I23').Value = '_PMC'

$ws.Range('E24').Value = '[Ya-Xi%Wang%NULL%1,     Hong-Tao%Guo%NULL%1,     Xue-Wei%Du%NULL%1,     Wen%Song%NULL%1,     Chang%Lu%NULL%1,     Wen-Nv%Hao%NULL%1,     Maria%Kapritsou.%NULL%3,     Maria%Kapritsou.%NULL%0,     Maria%Kapritsou.%NULL%0]'
$ws.Range('I24').Value = '_PMC'

$ws.Range('E25').Value = '[Yufei%Li%NULL%0,     Nathaniel%Scherer%NULL%0,     Lambert%Felix%NULL%0,     Lambert%Felix%NULL%0,     Hannah%Kuper%NULL%0,     Hannah%Kuper%NULL%0,     Jakob%Pietschnig%NULL%0,     Jakob%Pietschnig%NULL%0]'
$ws.Range('I25').Value = '_PMC'

$ws.Range('E26').Value = '[Xiao%Xiao%NULL%0,     Xiaobin%Zhu%NULL%1,     Shuai%Fu%NULL%1,     Yugang%Hu%NULL%1,     Xiaoning%Li%NULL%0,     Jinsong%Xiao%NULL%1]'
$ws.Range('I26').Value = '_PMC_elsevier'

$ws.Range('E27').Value = '[Stephen X.%Zhang%NULL%0,     Jing%Liu%NULL%0,     Asghar%Afshar Jahanshahi%NULL%1,     Khaled%Nawaser%NULL%1,     Ali%Yousefi%NULL%1,     Jizhen%Li%NULL%1,     Shuhua%Sun%NULL%1]'
$ws.Range('I27').Value = '_PMC_elsevier'

$ws.Range('E28').Value = '[Wen-rui%Zhang%NULL%0,     Kun%Wang%NULL%2,     Lu%Yin%NULL%2,     Wen-feng%Zhao%NULL%1,     Qing%Xue%NULL%1,     Mao%Peng%NULL%1,     Bao-quan%Min%NULL%1,     Qing%Tian%NULL%1,     Hai-xia%Leng%NULL%1,     Jia-lin%Du%NULL%1,     Hong%Chang%NULL%1,     Yuan%Yang%NULL%2,     Wei%Li%NULL%2,     Fang-fang%Shangguan%NULL%1,     Tian-yi%Yan%NULL%1,     Hui-qing%Dong%NULL%1,     Ying%Han%NULL%1,     Yu-ping%Wang%NULL%1,     Fiammetta%Cosci%NULL%1,     Hong-xing%Wang%NULL%1]'
$ws.Range('I28').Value = '_PMC'

$ws.Range('E29').Value = '[Yeen%Huang%NULL%0,     Ning%Zhao%zhaoning2018@email.szu.edu.cn%1]'
$ws.Range('I29').Value = '_PMC_elsevier'

$ws.Range('E30').Value = '[Juhong%Zhu%NULL%0,     Lin%Sun%NULL%1,     Lan%Zhang%NULL%1,     Huan%Wang%NULL%1,     Ajiao%Fan%NULL%1,     Bin%Yang%NULL%1,     Wei%Li%NULL%0,     Shifu%Xiao%NULL%1]'
$ws.Range('I30').Value = '_PMC'

$ws.Range('E31').Value = '[Jiang%Du%NULL%0,     Lu%Dong%NULL%1,     Tao%Wang%NULL%10,     Chenxin%Yuan%NULL%1,     Rao%Fu%NULL%1,     Lei%Zhang%NULL%2,     Bo%Liu%NULL%0,     Mingmin%Zhang%NULL%1,     Yuanyuan%Yin%NULL%1,     Jiawen%Qin%NULL%1,     Jennifer%Bouey%NULL%1,     Min%Zhao%NULL%1,     Xin%Li%NULL%2]'
$ws.Range('I31').Value = '_PMC_elsevier'

$ws.Range('E32').Value = '[Rümeysa Yeni%Elbay%rumeysa.yenielbay@medeniyet.edu.tr%0,     Ayşe%Kurtulmuş%NULL%1,     Selim%Arpacıoğlu%NULL%1,     Emrah%Karadere%NULL%1]'
$ws.Range('I32').Value = '_PMC_elsevier'

$ws.Range('E33').Value = '[Anne%Moorhead%NULL%1,     Gunther%Eysenbach%NULL%1,     Albert%Wu%NULL%2,     Albert%Wu%NULL%0,     Bojana%Knezevic%NULL%1,     Bradley A%Evanoff%bevanoff@wustl.edu%2,     Bradley A%Evanoff%bevanoff@wustl.edu%0,     Jaime R%Strickland%NULL%2,     Jaime R%Strickland%NULL%0,     Ann Marie%Dale%NULL%2,     Ann Marie%Dale%NULL%0,     Lisa%Hayibor%NULL%2,     Lisa%Hayibor%NULL%0,     Emily%Page%NULL%2,     Emily%Page%NULL%0,     Jennifer G%Duncan%NULL%2,     Jennifer G%Duncan%NULL%0,     Thomas%Kannampallil%NULL%2,     Thomas%Kannampallil%NULL%0,     Diana L%Gray%NULL%2,     Diana L%Gray%NULL%0]'
$ws.Range('I33').Value = '_PMC'

$ws.Range('E34').Value = '[Yufei%Li%NULL%0,     Nathaniel%Scherer%NULL%0,     Lambert%Felix%NULL%0,     Lambert%Felix%NULL%0,     Hannah%Kuper%NULL%0,     Hannah%Kuper%NULL%0,     Jakob%Pietschnig%NULL%0,     Jakob%Pietschnig%NULL%0]'
$ws.Range('I34').Value = '_PMC'

$ws.Range('E35').Value = '[Seshadri Sekhar%Chatterjee%NULL%0,     Ranjan%Bhattacharyya%NULL%1,     Sumita%Bhattacharyya%NULL%1,     Sukanya%Gupta%NULL%1,     Soumitra%Das%NULL%1,     Bejoy Bikram%Banerjee%NULL%1]'
$ws.Range('I35').Value = '_PMC'

$ws.Range('E36').Value = '[Jie%Chen%NULL%0,     Xinghuang%Liu%NULL%2,     Xinghuang%Liu%NULL%0,     Dongke%Wang%NULL%2,     Dongke%Wang%NULL%0,     Yan%Jin%NULL%2,     Yan%Jin%NULL%0,     Miao%He%NULL%1,     Yanling%Ma%NULL%1,     Xiaolong%Zhao%NULL%1,     Shuangning%Song%NULL%1,     Lei%Zhang%NULL%0,     Xuelian%Xiang%NULL%1,     Ling%Yang%NULL%1,     Jun%Song%song111jun@126.com%1,     Tao%Bai%drbaitao@126.com%4,     Xiaohua%Hou%houxh@hust.edu.cn%2,     Xiaohua%Hou%houxh@hust.edu.cn%0]'
$ws.Range('I36').Value = '_PMC_Springer'

$ws.Range('E37').Value = '[Yun%Chen%NULL%0,     Hao%Zhou%haoye320@163.com%1,     Yan%Zhou%NULL%1,     Fang%Zhou%NULL%1]'
$ws.Range('I37').Value = '_PMC_elsevier'

$ws.Range('I38').Value = ''

$ws.Range('E39').Value = '[Anucha%Apisarnthanarak%NULL%1,     Piyaporn%Apisarnthanarak%NULL%2,     Piyaporn%Apisarnthanarak%NULL%0,     Chanida%Siripraparat%NULL%2,     Chanida%Siripraparat%NULL%0,     Pavarat%Saengaram%NULL%1,     Narakorn%Leeprechanon%NULL%1,     David J.%Weber%NULL%1]'
$ws.Range('I39').Value = '_PMC'

$ws.Range('E40').Value = '[Andrea%Amerio%NULL%1,     Davide%Bianchi%NULL%1,     Francesca%Santi%NULL%1,     Luigi%Costantini%NULL%1,     Anna%Odone%NULL%1,     Carlo%Signorelli%NULL%1,     Alessandra%Costanza%NULL%1,     Gianluca%Serafini%NULL%1,     Mario%Amore%NULL%1,     Andrea%Aguglia%NULL%1]'
$ws.Range('I40').Value = '_PMC'

$ws.Range('E41').Value = '[Yufei%Li%NULL%0,     Nathaniel%Scherer%NULL%0,     Lambert%Felix%NULL%0,     Lambert%Felix%NULL%0,     Hannah%Kuper%NULL%0,     Hannah%Kuper%NULL%0,     Jakob%Pietschnig%NULL%0,     Jakob%Pietschnig%NULL%0]'
$ws.Range('I41').Value = '_PMC'
